# bts7960_wiring_diagram.pptx edit
#   - Footer "date" field on the slide master + every slide layout:
#       07-Jun-22 -> 23-Jun-22
#   - Title textbox: split the firmware-exception sentence into 3 runs and
#     rename fw-v183 -> fw-vXX3
#   - "BTS7960 ..." caption textbox: merge the two runs in each paragraph
#   - "If your motor is very hard..." textbox: re-flow the wording and fix
#     the "bts" -> "BTS" spelling (dropping the spell-check err flag)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text: master + all custom layouts
# ---------------------------------------------------------------------
function Set-DateText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newText
        }
    }
}

$master = $p.SlideMaster
Set-DateText $master.Shapes "23-Jun-22"

for ($l = 1; $l -le $master.CustomLayouts.Count; $l++) {
    $layout = $master.CustomLayouts.Item($l)
    Set-DateText $layout.Shapes "23-Jun-22"
}

# ---------------------------------------------------------------------
# 2) Title textbox ("Wiring diagram for H-bridge motor driver (valid ...)")
# ---------------------------------------------------------------------
$slide = $p.Slides.Item(1)

$titleBox = $slide.Shapes.Item("TextBox 107")
$titleRange = $titleBox.TextFrame.TextRange

# merge the old 2 runs ("H-bridge motor driver " + "(valid ... fw-v183)")
# into one run holding the final text, then re-split it into three runs.
$merged = $titleRange.Characters(20, 75)
$merged.Text = "H-bridge motor driver (valid for all firmware in PWM" + [char]0x00B1 + " mode, except fw-vXX3)"

# split off the trailing ")" first so offsets for the earlier split stay put
$closeParen = $titleRange.Characters(94, 1)
$closeParen.Text = ")"

# split off "fw-vXX3"
$fwRun = $titleRange.Characters(87, 7)
$fwRun.Text = "fw-vXX3"

# ---------------------------------------------------------------------
# 3) "BTS7960 or any other 2 channel H-bridge / motor driver with PWM
#    inputs" caption: merge the 2 runs in each of the 2 paragraphs.
# ---------------------------------------------------------------------
$capBox = $slide.Shapes.Item("TextBox 170")
$capRange = $capBox.TextFrame.TextRange

$cap1 = $capRange.Characters(1, 39)
$cap1.Text = "BTS7960 or any other 2 channel H-bridge"

$cap2 = $capRange.Characters(41, 28)
$cap2.Text = "motor driver with PWM inputs"

# ---------------------------------------------------------------------
# 4) "If your motor is very hard to turn..." textbox: re-flow 3
#    paragraphs and fix "bts" -> "BTS" (dropping the the err="1" flag).
# ---------------------------------------------------------------------
$noteBox = $slide.Shapes.Item("TextBox 128")
$noteRange = $noteBox.TextFrame.TextRange

# Para: "signal input, " + "then you may"  ->  one run
$para1 = $noteRange.Characters(65, 26)
$para1.Text = "signal input, then you may"

# Para: "try the wiring trick by " + "shorting"  ->  one run
$para2 = $noteRange.Characters(92, 32)
$para2.Text = "try the wiring trick by shorting"

# Para: "bts" (err=1) + " pins " + "VCC, R_EN, L_EN" -> "BTS " + "pins " + "VCC, R_EN, L_EN"
# Drop the misspelled run entirely (clears the err flag), then rebuild the
# "BTS " / "pins " split from the remaining (clean) run.
$btsRun = $noteRange.Characters(125, 3)
$btsRun.Text = ""

$pinsRun = $noteRange.Characters(125, 6)
$pinsRun.Text = "BTS pins "

$pinsTail = $noteRange.Characters(129, 5)
$pinsTail.Text = "pins "

Write-Host "edit complete"
